# Commitando parte de RG
# Adds a new "rg" worksheet (Repercussao Geral) right before the hidden
# "recebidos_classe (2)" sheet, and fills it with the historical RG table.

$wb = $excel.ActiveWorkbook

# --- Insert the new sheet in the right position -------------------------
$beforeSheet = $wb.Worksheets.Item("recebidos_classe (2)")
$ws = $wb.Worksheets.Add($beforeSheet)
$ws.Name = "rg"

# --- Header row -----------------------------------------------------------
$ws.Range("A1").Value = "RG"
$ws.Range("B1").Value = 2016
$ws.Range("C1").Value = 2017
$ws.Range("D1").Value = 2018
$ws.Range("E1").Value = 2019
$ws.Range("F1").Value = 2020

# --- Data rows --------------------------------------------------------------
$ws.Range("A2").Value = "Repercussão Geral Reconhecida"
$ws.Range("B2").Value = 27
$ws.Range("C2").Value = 39
$ws.Range("D2").Value = 32
$ws.Range("E2").Value = 41
$ws.Range("F2").Value = 25

$ws.Range("A3").Value = "Repercussão Geral Negada"
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 18

$ws.Range("A4").Value = "Mérito Julgado"
$ws.Range("B4").Value = 28
$ws.Range("C4").Value = 38
$ws.Range("D4").Value = 23
$ws.Range("E4").Value = 24
$ws.Range("F4").Value = 116

$ws.Range("A5").Value = "Reafirmação de Jurisprudência"
$ws.Range("B5").Value = 9
$ws.Range("C5").Value = 9
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 8

# --- Cosmetics: black font for the whole table + autofit column A ---------
$tableRange = $ws.Range("A1:F5")
$tableRange.Font.Color = 0
$ws.Columns.Item(1).AutoFit() | Out-Null

# --- Selection matches the source file (cursor parked on F4) --------------
$ws.Range("F4").Select() | Out-Null
